# "all test suite excel updated"
# Rewrite the TestSuite sheet's data table (A1:C15) with the refreshed
# set of test-suite rows, widen column A, bump the height of a few rows,
# and leave the selection where the author last left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("TSID",                        "Description",                                       "Runmode"),
    @("Login_Verification",          "All type of login execute",                         "Y"),
    @("Navigate_Verification",       "Navigation based Test Cases execute",               "N"),
    @("Product_Verification",        "All Products based Test Cases execute",             "N"),
    @("AddToCart_Verification",      "Add all type of product based Test Cases execute",  "N"),
    @("Order_Module",                "Order status  based Test Cases execute",            "N"),
    @("EvolveMoney_Module",          "Evolve money based Test Cases execute",             "N"),
    @("Referral _Module",            "Referal based Test Cases execute",                  "N"),
    @("GiftVoucher_Module",          "Gift Voucher based Test Cases",                     "N"),
    @("Hamper_Module",               "Hamper_Module based Test Cases",                    "N"),
    @("Subscription_Module",         "Subscription_Module based Test Cases",              "N"),
    @("Form_Verification",           "Form_Verification based Test Cases",                "N"),
    @("Link_Verification",           "Link_Verification based Test Cases",                "N"),
    @("Browser_Verification",        "Browser_Verification based Test Cases",             "N"),
    @("Registeration_Verification",  "Registeration_Verification based Test Cases",       "N")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Column A needs to be noticeably wider to fit the new, longer module names.
$ws.Columns.Item(1).ColumnWidth = 36.6

# A few rows (7-9) now carry an explicit row height.
$ws.Rows.Item(7).RowHeight = 15.75
$ws.Rows.Item(8).RowHeight = 15.75
$ws.Rows.Item(9).RowHeight = 15.75

# Leave the cursor where the author left it after editing.
$ws.Range("B18").Select() | Out-Null
